# Agregar Codificacion del tipo de cuenta
# Inserts a new "AcctNumberType_EnumType" sheet as the first sheet of the
# workbook and populates it with the OECD account-number-type coding table
# (COD / DESC / OBSER columns).

$wb = $excel.ActiveWorkbook

# --- incidental state left on another sheet before the new sheet was added ---
# (the workbook previously had the cursor resting on CrsPaymentType_EnumType!B5)
$crsPayment = $wb.Worksheets.Item("CrsPaymentType_EnumType")
$crsPayment.Activate()
$crsPayment.Range("B5").Select()

# --- add the new sheet as the very first tab ---
$ws = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$ws.Name = "AcctNumberType_EnumType"

# --- headers ---
$ws.Range("A1").Value = "COD"
$ws.Range("B1").Value = "DESC"

# --- column A: codes ---
$ws.Range("A2").Value = "OECD601"
$ws.Range("A3").Value = "OECD602"
$ws.Range("A4").Value = "OECD603"
$ws.Range("A5").Value = "OECD604"
$ws.Range("A6").Value = "OECD605"

# --- column B: short type ---
$ws.Range("B2").Value = "IBAN"
$ws.Range("B3").Value = "OBAN"
$ws.Range("B4").Value = "ISIN"
$ws.Range("B5").Value = "OSIN"
$ws.Range("B6").Value = "Other"

# --- column C: description ---
$ws.Range("C2").Value = "International Bank Account Number"
$ws.Range("C3").Value = "Other Bank Account Number"
$ws.Range("C4").Value = "International Securities Information Number"
$ws.Range("C5").Value = "Other Securities Information Number"
$ws.Range("C6").Value = "Any Other type of account number"

# --- column C header, added last ---
$ws.Range("C1").Value = "OBSER"

# widen column C to fit the long description text
$ws.Columns.Item(3).ColumnWidth = 40.6

# leave the new sheet active/selected
$ws.Activate()
$ws.Range("A1").Select()
